# Update the COVID-19 "Pais" (countries) dashboard with the next data refresh.
#
# The refresh brings in new daily totals for several countries. Because the
# sheet is kept sorted by "Casos totales" (column B, descending), a handful
# of countries that overtake their neighbours in total cases swap rows with
# them:
#   - Catar overtakes Panama            (rows 46/47)
#   - Kazajistan overtakes Uzbekistan   (rows 69/70)
#   - Sierra Leona overtakes Botsuana / San Vicente y las Granadinas / Seychelles
#     and jumps from row 192 up to row 189
#
# Brasil (row 17) and Portugal (row 19) simply get fresh numbers without
# changing rank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a full data row (country name + the 7 metric columns).
# NOTE: positional parameters are used throughout because this runtime does
# not bind named (-Param value) arguments.
function Set-CountryRow {
    param($Row, $Country, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes)

    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Header: refresh the "last updated" timestamp -------------------------
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 13:52"

# --- Plain data refresh (no re-ordering) -----------------------------------
Set-CountryRow 17 "Brasil"   25758 496 14026 10175 296 25 1557
Set-CountryRow 19 "Portugal" 18091 643 383   17109 208 32 599

# --- Catar overtakes Panama (rows 46/47) ------------------------------------
Set-CountryRow 46 "Catar"  3711 283 406 3298 37  0 7
Set-CountryRow 47 "Panama" 3574 0   72  3407 106 0 95

# --- Kazajistan overtakes Uzbekistan (rows 69/70) ---------------------------
Set-CountryRow 69 "Kazajistan" 1290 58  220 1055 20 1 15
Set-CountryRow 70 "Uzbekistan" 1275 110 99  1172 8  0 4

# --- Sierra Leona jumps from row 192 to row 189 -----------------------------
# Botsuana and San Vicente y las Granadinas keep their own (unchanged) data,
# they just shift down one row; Seychelles ends up with the data that used
# to belong to Sierra Leona's old row (which happens to already match).
Set-CountryRow 189 "Sierra Leona"                 13 2 0 13 0 0 0
Set-CountryRow 190 "Botsuana"                     13 0 0 12 0 0 1
Set-CountryRow 191 "San Vicente y las Granadinas" 12 0 1 11 0 0 0
Set-CountryRow 192 "Seychelles"                   11 0 0 11 0 0 0
